$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 6271.3313874533
$ws.Range("C2").Value = 12657.0639087956
$ws.Range("D2").Value = 28519.37911664036
